# Normalize the "Recorded By" (column G) entries on the active sheet.
# For every data row, if the cell contains a comma-separated list of
# recorder names/emails whose first entry is not exactly (case-sensitive)
# "System", rotate the list left by one position (move the first entry
# to the end). Rows that already start with "System" (or contain a
# single value) are left untouched.

function Test-CaseSensitiveEquals($s1, $s2) {
    if ($s1.Length -ne $s2.Length) { return $false }
    for ($i = 0; $i -lt $s1.Length; $i++) {
        if ([int][char]$s1[$i] -ne [int][char]$s2[$i]) {
            return $false
        }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $cur = $cell.Value2

    if ($cur -ne $null -and $cur -ne "") {
        $parts = $cur -split ", "

        if ($parts.Length -gt 1 -and -not (Test-CaseSensitiveEquals $parts[0] "System")) {
            $rotated = $parts[1..($parts.Length - 1)] + $parts[0]
            $newVal = $rotated -join ", "

            if ($newVal -ne $cur) {
                $ws.Cells.Item($r, 7).Value = $newVal
            }
        }
    }
}
